$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 7.217020034790039
$ws.Range("B1").Value = 5.770425796508789
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 4.544404029846191
$ws.Range("E1").Value = 2.114849805831909
